$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-04"

# Update the November row label to reflect the new "through" date
$ws.Range("A12").Value = "November (through 11-04)"

# Update November row (row 12) values for each year column (B:I)
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 16
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 28
$ws.Range("H12").Value = 24
$ws.Range("I12").Value = 13

# Update Total row (row 13) values for each year column (B:I)
$ws.Range("B13").Value = 261
$ws.Range("C13").Value = 496
$ws.Range("D13").Value = 726
$ws.Range("E13").Value = 631
$ws.Range("F13").Value = 488
$ws.Range("G13").Value = 1085
$ws.Range("H13").Value = 1465
$ws.Range("I13").Value = 1413
